# Hesperia.xlsx - "Master Scores" sheet update
# Commit message: "Python code completed (Brute force)"
#
# Changes applied:
#   1. B2 and C2 become -1000 (brute-force search sentinel / penalty values)
#   2. F16 becomes 40 (was 43)
#   3. The saved selection moves to C3 (was E25)
#   4. Cosmetic: the workbook window geometry is nudged (best effort - the
#      headless runtime's ActiveWindow geometry is view-only state and may
#      not be persisted to the saved file, but we still set it through the
#      documented object model in case it is honored).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data edits -----------------------------------------------------------
$ws.Range("B2").Value = -1000
$ws.Range("C2").Value = -1000
$ws.Range("F16").Value = 40

# --- Selection --------------------------------------------------------------
$ws.Range("C3").Select()

# --- Window geometry (best effort) ------------------------------------------
$win = $excel.ActiveWindow
$win.Left = 14820
$win.Top = 0
$win.Width = 14000
$win.Height = 17460

$wb.Save()
